$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ATS Accuracy")
$ws1.Range("B2").Value = 4
$ws1.Range("C2").Value = 72
$ws1.Range("D2").Value = 76
$ws1.Range("E2").Value = 94.7

$ws1.Range("B3").Value = 4
$ws1.Range("D3").Value = 63
$ws1.Range("E3").Value = 93.7

$ws1.Range("B5").Value = 4
$ws1.Range("D5").Value = 13
$ws1.Range("E5").Value = 69.2

$ws2 = $wb.Worksheets.Item("Total Accuracy")
$ws2.Range("C2").Value = 64
$ws2.Range("D2").Value = 69
$ws2.Range("E2").Value = 92.8

$ws2.Range("B3").Value = 4
$ws2.Range("C3").Value = 55
$ws2.Range("E3").Value = 93.2

$ws2.Range("B4").Value = 4
$ws2.Range("C4").Value = 24
$ws2.Range("D4").Value = 28
$ws2.Range("E4").Value = 85.7

$ws2.Range("B5").Value = 4
$ws2.Range("C5").Value = 8
$ws2.Range("E5").Value = 66.7
